# CROX_Model.xlsx edit: swap the "EPS" row (31) and the "Shares Outstanding
# (SEC 1st page)" row (32) on the Model sheet, fixing up the EPS formula to
# reference the new location of the shares-outstanding row, and switch the
# active sheet/selection over to the Model sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")

# ---------------------------------------------------------------------
# 1) Row 31 / Row 32 content swap
#    Before: row31 = "EPS" (=rowN30/rowN32 shared formula, columns C:T)
#            row32 = "Shares Outstanding (SEC 1st page)" (literal values)
#    After:  row31 = "Shares Outstanding (SEC 1st page)" (literal values)
#            row32 = "EPS" (=rowN30/rowN31 formula, columns C:T)
# ---------------------------------------------------------------------

# Shares-outstanding literal values, captured from the original row 32.
$sharesOutstanding = [ordered]@{
    "H" = 62.386000000000003
    "I" = 58.847000000000001
    "K" = 61.58
    "L" = 61.65
    "M" = 61.744999999999997
    "N" = 61.750999999999998
    "O" = 62.026000000000003
    "P" = 59.384999999999998
    "Q" = 60.567
    "R" = 60.499000000000002
    "S" = 60.703000000000003
    "T" = 59.386000000000003
}

# --- Write "Shares Outstanding (SEC 1st page)" into row 31 ---
$ws.Range("B31").Value = "Shares Outstanding (SEC 1st page)"
$ws.Range("C31:G31").ClearContents()
$ws.Range("J31").ClearContents()
foreach ($col in $sharesOutstanding.Keys) {
    $ws.Range(($col + "31")).Value = $sharesOutstanding[$col]
}

# --- Write "EPS" formula row into row 32 (denominator now row 31) ---
# Each formula is entered individually (not as one range fill) so every
# cell gets its own <f> element rather than a shared-formula group.
$ws.Range("B32").Value = "EPS"
foreach ($col in @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")) {
    $ws.Range(($col + "32")).Formula = "=" + $col + "30/" + $col + "31"
}

# --- Fix up the blank formatting cells U/V that travelled with the swap ---
$ws.Range("V31").Style = "Normal"
$ws.Range("V31").Font.Name = "Arial"
$ws.Range("V31").NumberFormat = "0.00"

$ws.Range("V32").Style = "Normal"
$ws.Range("V32").Font.Name = "Arial"

# ---------------------------------------------------------------------
# 2) View state: make "Model" the active sheet/tab, restore the frozen
#    pane split (2,2) and move the selection to T33.
# ---------------------------------------------------------------------
$ws.Range("C3").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("T33").Select()
